# Adds a new header/annotation row (row 1) above the existing field-name
# row (row 2) of the device-inventory upload template, with per-column
# notes (required/optional, allowed enum values, defaults) and turns on
# word-wrap for the note cells. Also nudges the active selection like the
# author's last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- reusable note strings -------------------------------------------------
$notReq      = "[not reqired]"
$hostname    = "Hostname`n[required,`nunique]"
$statusEnum  = "[In operation,`nDecommissioned,`nPending Setup,`nOffline,`nNot defined yet,`nException]"
$categoryEnum = "[Network,`nServer - Virtual,`nServer - Physical,`n Server - Platform,`nPrinter/Scanner, `nStorage, `nConferencing, `nEnd Users Computing, `nAirCon, `nUPS, `nOther]"
$subCategoryEnum = "[Router,`nFirewall,`nFirewall/IDS/IPS,`nAccess Point,`nSwitch,`nDesktop,`nLaptop,`nPrinter/MFP,`nScanner,`nIP Phone,`nTeleconferencing/Modem,`nVoIP System - Cisco CM,`nVoIP System - other,`nApp and DB Server,`nApplication Server,`nDatabase Server,`nFile Server,`nOther Server,`nBackup device,`nStorage - NAS,`nStorage - SAN,`nData Historian,`nHuman Machine Interface (HMI),`nIDS/IPS Detection,`nMaster Terminal Unit (MTU),`nProgrammable Logic Controller (PLC),`nRemote Access,`nRemote Terminal Unit (RTU),`nOther hardware]"
$default1    = "default=1`n[not reqired]"
$defaultToday = "default=today`n[not reqired]"

# ---- row 1 values ----
# Note: new distinct strings are appended to the shared-strings table in
# first-write order, so the write order below (not just the final layout)
# is chosen to reproduce the original author's shared-string indices:
# Hostname, [not reqired], status-enum, category-enum, sub_category-enum,
# default=1, default=today.
$ws.Range("A1").Value = $hostname
$ws.Range("B1").Value = $notReq
$ws.Range("C1").Value = $notReq
$ws.Range("D1").Value = $statusEnum
$ws.Range("E1").Value = $categoryEnum
$ws.Range("F1").Value = $subCategoryEnum
$ws.Range("G1").Value = $notReq
$ws.Range("H1").Value = $notReq
$ws.Range("I1").Value = $notReq
$ws.Range("J1").Value = $notReq
$ws.Range("K1").Value = $notReq
$ws.Range("L1").Value = $notReq
$ws.Range("M1").Value = $notReq
$ws.Range("N1").Value = $notReq
$ws.Range("O1").Value = $notReq
$ws.Range("P1").Value = $notReq
$ws.Range("Q1").Value = $notReq
$ws.Range("V1").Value = $default1
$ws.Range("W1").Value = $default1
$ws.Range("R1").Value = $defaultToday
$ws.Range("S1").Value = $defaultToday
$ws.Range("T1").Value = $defaultToday
$ws.Range("U1").Value = $notReq

# ---- wrap text + row height for the new note row ---------------------------
# (applied cell-by-cell: a multi-area Range's .WrapText setter only takes
# effect on the first area in this engine)
foreach ($addr in @("A1","D1","E1","F1","R1","S1","T1","V1","W1")) {
    $ws.Range($addr).WrapText = $true
}
$ws.Rows.Item(1).RowHeight = 409.5

# ---- column widths (approximate the widened columns from the author's
# AutoFit pass; the engine snaps widths to its own pixel grid so we pick
# the nearest supported ColumnWidth to each target) -------------------------
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 19
$ws.Columns.Item(5).ColumnWidth = 15
$ws.Columns.Item(6).ColumnWidth = 30
$ws.Columns.Item(18).ColumnWidth = 13
$ws.Columns.Item(19).ColumnWidth = 13
$ws.Columns.Item(20).ColumnWidth = 13
$ws.Columns.Item(22).ColumnWidth = 12

# ---- selection, as last left by the author ---------------------------------
$ws.Range("J9").Select()
